$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 - this shifts the existing rows 2-5 down to 3-6.
$ws.Rows.Item(2).Insert()

# Fill in the new "Terapeuta" job row (row 2).
$ws.Range("A2").Value = 'Terapeuta'
$ws.Range("B2").Value = 'La empresa es confidencial o no se encuentra disponible'
$ws.Range("C2").Value = 'León,, Gto.'
$ws.Range("D2").Value = '$14,000 - $16,000 Mensual'
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = 'Sector salud'
$ws.Range("G2").Value = 'Terapeuta'
$ws.Range("H2").Value = 'Universitario titulado'
$ws.Range("I2").Value = 'Permanente'
$ws.Range("J2").Value = 'Tiempo completo'
$ws.Range("K2").Value = 'Presencial'

$descripcion = @'
Requisitos del puesto
Estudios universitarios con título en Terapia.
Experiencia previa como Terapeuta de niños con Trastornos del espectro autista.
Gusto por realizar manualidades.
Habilidad para nadar.
Licencia de manejo vigente.
Conocimientos en técnicas de terapia y rehabilitación.
Licencia o certificación válida en Terapia (deseable).
Responsabilidades del puesto
Realizar evaluaciones y diagnósticos de los pacientes.
Diseñar planes de tratamiento personalizados.
Realizar sesiones de terapia adaptadas a las necesidades individuales de cada paciente.
Mantener registros precisos de la evolución de los pacientes.
Prestaciones y beneficios adicionales
Salario mensual competitivo de 14000 a 16000.
Prestaciones de ley.
Vales de despensa.
Fondo de ahorro.
Oportunidades de capacitación y desarrollo profesional en un ambiente de trabajo colaborativo y en constante crecimiento.
'@
$ws.Range("L2").Value = $descripcion

# After the insert, the original rows are:
#  row 3 = Maestra (unchanged)
#  row 4 = Monitora en inclusión educativa (maestro sombra) (unchanged)
#  row 5 = TERAPEUTA ESPECIALIZADA EN TEA Y TDAH (to be removed)
#  row 6 = Closer SaaS B2B Startup Tecnologica (to be removed)
$ws.Range("A5:A6").EntireRow.Delete()
